$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Fields/FormOID): replace Rules (C3) and Questions (D3); add Name of DB function (E3) ---
$c3 = @'
Use the caDSR Form Long Name that the Questions are associated with, in ALL CAPS; Should be the same name a in the OID column on Forms Tab. 
'@
$ws.Range("C3").Value = $c3
$d3 = @'
If CDEs (not forms) are being imported, then "Field = CDE". Correct?

'@
$ws.Range("D3").Value = $d3
$ws.Range("E3").Value = "FUNCTION ONEDATA_WA.ALS_LONG_NAME"

# --- Insert new row 4 (FieldOID) ---
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Clear()
$ws.Range("B4").Value = "FieldOID"
$c4 = @'
If Forms are downloaded: FORM_OID 
If CDEs are downloaded: UserID + "_CDECART"
One row for each question on the form UNLESS the Question is a "all that apply" type Question.
For questions with "all" in instruction:
Create the first field for the question with FieldOID as below, + "_An" where n a number starting with "1" + _LBLAnn" where nn = n e.g. "PT_RACE_CD_A1_LBLA1" 
Create the field for each value in the Question with FieldOID  = same rule as below + "_An" where "n" is number starting with 1 for the first value, incremented by 1 for each value. E.g. If there are 5 values, _A1, _A2, _A3, _A4, _A5
e.g. PT_RACE_CD_A1, PT_RACE_CD_A2, PT_RACE_CD_A3, PT_RACE_CD_A4, PT_RACE_CD_A5
FieldOID Must be unique within the all the form.
Limited to 30 characters.
Use the caDSR CDE Short Name that is linked to the Question. Use the Question-CDE link to retreive the CDE short name. 
 If no link to a CDE then default is UserName + "_" + integer (start with 1). 
When Questions are used multiple times on a caDSR form, must create these questions as unique fields by adding a sequence number to the  CDE Short Name e.g. “RACE_01”, to ensure uniqueness in Rave.This ensures uniqueness Field OIDs 
Replace spaces and other characters:
  ~ Only underscore characters are permitted
  ~ Replace space characters with “_”
  ~ Replace “.” with “_”
  ~ Replace “/” with “_”
  ~ Replace “(” with “_”
  ~ Replace “)” with “_”
  ~ Replace “-” with “_”
  ~ Replace “’” with “”
'@
$ws.Range("C4").Value = $c4
$ws.Range("D4").Value = "FORM_OID or Question ID?"
$ws.Rows.Item(4).RowHeight = 409.5

# --- Row 3 height update ---
$ws.Rows.Item(3).RowHeight = 43.5

# --- Row 6 (was row5): DraftFieldName -- add Long name question (D6) ---
$ws.Range("D6").Value = "Long name can be up to 255. How to treat then concatination?"

# --- Row 7 (was row6): VariableOID -- add FormOID/FORM_OID question (D7) ---
$ws.Range("D7").Value = "What is the difference between FormOID and FORM_OID?"

# --- New rows 8 and 9 ---
$ws.Range("B8").Value = "DataFormat"
$ws.Range("C8").Value = "?"
$ws.Range("B9").Value = "DataDictionaryName"

# --- Column B width ---
$ws.Columns.Item(2).ColumnWidth = 20.7265625

# --- View: active cell D4 ---
$ws.Range("D4").Select()
